$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first land record (04880023, index 14).
# This shifts the remaining record (04880090, index 15) up from row 3 to row 2,
# and the engine will prune now-unused shared strings automatically.
$ws.Rows(2).Delete()

# Add two new trailing columns: "portion" and "total" (area * share portion) of land.
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("P1").Value = "portion"

$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("Q1").Value = "total"

$ws.Range("C2").Copy($ws.Range("P2"))
$ws.Range("P2").Value = 1

$ws.Range("C2").Copy($ws.Range("Q2"))
$ws.Range("Q2").Value = 1
